$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: consequents_length, mirroring the antecedents_length header (G)
# Copy G1's formatting (bold header style) onto H1, then set its value/label.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "consequents_length"

# Data rows for the new column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
